$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Add a new "description" header in column M, row 1 (next to existing "b" header in L1)
$ws.Range("M1").Value = "description"

# Move/collapse the active selection to the newly added header cell
$ws.Range("M1").Select()
